$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.916.04'
$ws.Range("E2").Value = '  -1.17%  '

$ws.Range("D3").Value = '1.628.99'
$ws.Range("E3").Value = '  -2.63%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("E4").Value = '  -0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.46'
$ws.Range("E5").Value = '  -0.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5235'
$ws.Range("E6").Value = '  -0.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2582'
$ws.Range("E8").Value = '  -2.46%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06280'
$ws.Range("E9").Value = '  -0.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.54'
$ws.Range("E10").Value = '  -3.76%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07571'
$ws.Range("E11").Value = '  +0.11%  '

$ws.Range("D12").Value = '1.629.34'
$ws.Range("E12").Value = '  -2.88%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.418'
$ws.Range("E13").Value = '  -0.94%  '

$ws.Range("D14").Value = '1.852.44'
$ws.Range("E14").Value = '  -2.55%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5523'
$ws.Range("E15").Value = '  -1.39%  '

$ws.Range("D16").Value = '0.0₅8026'
$ws.Range("E16").Value = '  -0.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.75'
$ws.Range("E17").Value = '  -3.47%  '

$ws.Range("D18").Value = '25.923.62'
$ws.Range("E18").Value = '  -1.27%  '

$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.676'
$ws.Range("E20").Value = '  -2.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '185.85'
$ws.Range("E21").Value = '  -0.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.14'
$ws.Range("E22").Value = '  -2.68%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.123'
$ws.Range("E23").Value = '  -1.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9993'
$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.14'
$ws.Range("E25").Value = '  -3.22%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1217'
$ws.Range("E26").Value = '  -3.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.383'
$ws.Range("E27").Value = '  -2.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.71'
$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.363'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05899'
$ws.Range("E30").Value = '  -4.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.247'
$ws.Range("E31").Value = '  -2.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.418'
$ws.Range("E32").Value = '  -2.49%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.389'
$ws.Range("E33").Value = '  -1.25%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.624'
$ws.Range("E34").Value = '  -0.63%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9796'
$ws.Range("E35").Value = '  -2.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.381'
$ws.Range("E36").Value = '  -1.30%  '

$ws.Range("E37").Value = '  -0.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5793'
$ws.Range("E38").Value = '  -4.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01602'
$ws.Range("E39").Value = '  -1.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8478'
$ws.Range("E40").Value = '  -3.21%  '

$ws.Range("D42").Value = '1.036.62'
$ws.Range("E42").Value = '  -5.98%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.667'
$ws.Range("E43").Value = '  -7.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.83'
$ws.Range("E44").Value = '  -0.14%  '

$ws.Range("D45").Value = '1.778.64'
$ws.Range("E45").Value = '  -2.43%  '

$ws.Range("D46").Value = '0.0₈106'
$ws.Range("E46").Value = '  -1.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.006'
$ws.Range("E47").Value = '  -0.03%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.11'
$ws.Range("E48").Value = '  -1.63%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.033'
$ws.Range("E49").Value = '  +0.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05165'
$ws.Range("E50").Value = '  -1.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4214'
$ws.Range("E51").Value = '  -0.96%  '
